$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.566.07"
$ws.Range("E2").Value = "  +2.66%  "
$ws.Range("D3").Value = "2.032.32"
$ws.Range("E3").Value = "  +7.25%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'245.25"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "'0.659"
$ws.Range("E6").Value = "  -4.94%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'44.50"
$ws.Range("E8").Value = "  +3.22%  "
$ws.Range("D9").Value = "'60.29"
$ws.Range("E9").Value = "  +5.84%  "
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").Value = "'0.0715"
$ws.Range("E11").Value = "  -5.11%  "
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "'14.35"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").Value = "2.326.64"
$ws.Range("E14").Value = "  +7.18%  "
$ws.Range("D15").Value = "'0.803"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").Value = "2.029.42"
$ws.Range("E16").Value = "  +7.03%  "
$ws.Range("D17").Value = "'4.86"
$ws.Range("E17").Value = "  -3.62%  "
$ws.Range("D18").Value = "36.514.32"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("D19").Value = "'70.87"
$ws.Range("E19").Value = "  -3.72%  "
$ws.Range("D20").Value = "0.0₃0809"
$ws.Range("E20").Value = "  -2.64%  "
$ws.Range("D21").Value = "'236.16"
$ws.Range("E21").Value = "  -4.10%  "
$ws.Range("D22").Value = "'12.55"
$ws.Range("E22").Value = "  -3.51%  "
$ws.Range("D23").Value = "'4.87"
$ws.Range("E23").Value = "  -6.46%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("E25").Value = "  -9.05%  "
$ws.Range("D26").Value = "'168.81"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'19.90"
$ws.Range("E27").Value = "  +8.26%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'8.70"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "'1.93"
$ws.Range("E29").Value = "  -9.94%  "
$ws.Range("E30").Value = "  -5.46%  "
$ws.Range("D31").Value = "'21.52"
$ws.Range("E31").Value = "  +51.15%  "
$ws.Range("D32").Value = "'4.33"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").Value = "'0.0577"
$ws.Range("E33").Value = "  -5.23%  "
$ws.Range("D34").Value = "'0.0897"
$ws.Range("E34").Value = "  +21.22%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("E36").Value = "  +1.41%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'3.96"
$ws.Range("E37").Value = "  -7.20%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.19"
$ws.Range("E38").Value = "  +12.01%  "
$ws.Range("D39").Value = "'0.851"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("E40").Value = "  -11.51%  "
$ws.Range("D41").Value = "'95.60"
$ws.Range("E41").Value = "  -3.72%  "
$ws.Range("D42").Value = "'0.0213"
$ws.Range("E42").Value = "  -7.28%  "
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("D44").Value = "'2.78"
$ws.Range("E44").Value = "  +15.61%  "
$ws.Range("D45").Value = "'15.73"
$ws.Range("E45").Value = "  -7.75%  "
$ws.Range("D46").Value = "1.311.43"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("D47").Value = "'0.0815"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("D49").Value = "2.227.47"
$ws.Range("E49").Value = "  +7.54%  "
$ws.Range("D50").Value = "'2.19"
$ws.Range("E50").Value = "  -7.18%  "
$ws.Range("D51").Value = "'3.79"
$ws.Range("E51").Value = "  +14.42%  "
